$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add "Done" status to a few rows that previously had no status, and
#     flip one row from "begonnen" to "done" (E33) -------------------------
$ws.Range("E6").Value = "done"
$ws.Range("E26").Value = "done"
$ws.Range("E27").Value = "done"
$ws.Range("E33").Value = "done"

# --- Move the active selection from E5 to C8 ------------------------------
[void]$ws.Range("C8").Select()

# --- Apply the AutoFilter on A2:E34, filtering column E (index 5, 1-based)
#     down to blanks + "begonnen" + "obsolet" (hides every "done" row) ----
$ws.Range("A2:E34").AutoFilter(5, @("begonnen", "obsolet", ""), 7)

# --- Hide every data row whose status is "done" (this is what Excel does
#     on-disk once that AutoFilter is applied) -----------------------------
$hiddenRows = @(6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 26, 27, 29, 30, 31, 33)
foreach ($r in $hiddenRows) {
    $ws.Rows.Item($r).Hidden = $true
}

# --- Register the hidden _FilterDatabase defined name that Excel creates
#     automatically for a sheet-level AutoFilter --------------------------
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$2:`$E`$34")
$fdb.Visible = $false
